# "massive MDY site template update"
# Insert three new columns (Month, Day, Year) before the existing
# "Date Sampled" column on the "Data Entry" sheet, and populate them from
# the Month/Day/Year parts of each row's existing sampling date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank columns at E:G (old E "Date Sampled" and everything to its
# right shifts right by 3, E->H, F->I, ... N->Q). The new columns pick up
# the same width as the "Transect" column (D) immediately to their left.
$ws.Columns("E:G").Insert()

# Header row.
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# Data rows - Month / Day / Year pulled from each record's sample date.
$ws.Range("E2:G5").Style = "Normal"

$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 2016

$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 11
$ws.Range("G3").Value = 2016

$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 2016

$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 2016

# Restore the active selection to the new Month column on row 5.
$ws.Range("E5").Select() | Out-Null
